$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (pushes existing rows 15..40 down to 16..41)
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new data record
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C15").Value = "Los Lagos"
$ws.Range("D15").Value = 44645
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 100112030
$ws.Range("G15").Value = "Poroto granado"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 70
$ws.Range("K15").Value = 26000
$ws.Range("L15").Value = 26000
$ws.Range("M15").Value = 26000
$ws.Range("N15").Value = "$/saco 25 kilos"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 1040
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
